# fix(gui) step 1 and 2
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Advance the quote date by one day (A1 holds a date serial)
$ws.Range("A1").Value = 45309

# Update "step 1" and "step 2" prices in column D
$ws.Range("D32").Value = 7320
$ws.Range("D33").Value = 8170
$ws.Range("D34").Value = 9280
$ws.Range("D35").Value = 11550
$ws.Range("D36").Value = 11960
$ws.Range("D37").Value = 12590
